$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "scutellarin only"

$ws.Range("A32").Value = "#"
$ws.Range("B32").Value = "mix"
$ws.Range("C32").Value = "chemical"
$ws.Range("D32").Value = "base conc (ppm)"
$ws.Range("E32").Value = "desired final vol (mL)"
$ws.Range("F32").Value = "to make 100 ppm mix (uL)"

$ws.Range("C33").Value = "scutellarin"
$ws.Range("D33").Value = 4000
$ws.Range("E33").Value = 1.5
$ws.Range("F33").Formula = "=`$E33*(100*0.001)/(`$D33*0.001)*1000"

$ws.Range("E34").Value = "50% MeOH to add (microL)"
$ws.Range("F34").Formula = "=(E33*1000)-SUM(F33)"

$ws.Range("A31").Font.Bold = $true
$ws.Range("A32:F32").Font.Bold = $true
$ws.Range("E34").Font.Bold = $true

$ws.Range("M16").Select() | Out-Null

$ws.PageSetup.Orientation = 1
